$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (IMEI shifts from B to C)
$ws.Columns("B").Insert()

# Set new header for inserted column B
$ws.Range("B1").Value = "LINEID"

# Match column B width to column A's custom width
$ws.Columns("B").ColumnWidth = 11

# Update selection to E3
$ws.Range("E3").Select()
